$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new "team record" columns (AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (e.g. A1) so the new
# headers match the bold/centered/bordered formatting used by the rest of
# row 1.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Fill in the win/loss/tie record for every data row (2-58). The 2019
# Philadelphia Phillies finished the season with an 81-81-0 record, so the
# same values are applied to each player row.
for ($row = 2; $row -le 58; $row++) {
    $ws.Cells.Item($row, 30).Value = 81
    $ws.Cells.Item($row, 31).Value = 81
    $ws.Cells.Item($row, 32).Value = 0
}
